# Insert a new daily price record for "Cebollín" at Terminal La Palmera de
# La Serena. This pushes the existing row 168 (and everything after it)
# down by one row, so we insert a fresh row at position 168 and populate
# it with the new observation. The sheet's used dimension grows from
# A1:R247 to A1:R248 automatically as part of the insert.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(168).Insert()

$ws.Cells.Item(168, 1).Value  = 8
$ws.Cells.Item(168, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(168, 3).Value  = "Coquimbo"
$ws.Cells.Item(168, 4).Value  = 44846
$ws.Cells.Item(168, 5).Value  = 4
$ws.Cells.Item(168, 6).Value  = 100112037
$ws.Cells.Item(168, 7).Value  = "Cebollín"
$ws.Cells.Item(168, 8).Value  = "Sin especificar"
$ws.Cells.Item(168, 9).Value  = "Primera"
$ws.Cells.Item(168, 10).Value = 1300
$ws.Cells.Item(168, 11).Value = 1400
$ws.Cells.Item(168, 12).Value = 1600
$ws.Cells.Item(168, 13).Value = 1500
$ws.Cells.Item(168, 14).Value = "$/paquete 6 unidades"
$ws.Cells.Item(168, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(168, 16).Value = 250
$ws.Cells.Item(168, 17).Value = 6
$ws.Cells.Item(168, 18).Value = "Hortaliza"
